$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.349.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6979"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.37%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3068"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07435"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.71"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08115"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.953.54"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7257"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.217"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.519.02"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.910"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007719"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.63%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.138.73"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.622"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1487"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.018"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.57"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.07"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.98%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.388"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.31%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.412"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.065"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05281"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.200"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7212"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.000"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.683"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01864"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.725"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8767"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4310"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.915"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.03"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.039.86"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.37"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.80%  "

$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.076.57"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.56%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.265"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.759"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.233"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.51%  "
